# Weekly update: insert a new daily price record for
# "Vega Monumental Concepción - Zanahoria" at row 413, shifting the
# existing historical rows (413-492) down by one (to 414-493).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 413, pushing rows 413:492 down to 414:493.
$ws.Rows.Item(413).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A413").Value = 11
$ws.Range("B413").Value = "Vega Monumental Concepción"
$ws.Range("C413").Value = "Bíobío"
$ws.Range("D413").Value = 45244
$ws.Range("E413").Value = 8
$ws.Range("F413").Value = 100114013
$ws.Range("G413").Value = "Zanahoria"
$ws.Range("H413").Value = "Sin especificar"
$ws.Range("I413").Value = "Primera"
$ws.Range("J413").Value = 250
$ws.Range("K413").Value = 5000
$ws.Range("L413").Value = 5500
$ws.Range("M413").Value = 5200
$ws.Range("N413").Value = "$/saco 20 kilos"
$ws.Range("O413").Value = "Región Metropolitana"
$ws.Range("P413").Value = 260
$ws.Range("Q413").Value = 20
$ws.Range("R413").Value = "Hortaliza"
